$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.757.22"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.855.71"
$ws.Range("E3").Value = "  +0.26%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.034"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'322.80"
$ws.Range("E5").Value = "  +0.42%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.07%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4399"
$ws.Range("E7").Value = "  +0.38%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3816"
$ws.Range("E8").Value = "  +1.68%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.07439"
$ws.Range("E9").Value = "  +0.43%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "'0.8869"
$ws.Range("E10").Value = "  +1.19%  "

# Row 11 - Solana
$ws.Range("D11").Value = "'21.58"
$ws.Range("E11").Value = "  +0.24%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.857.94"
$ws.Range("E12").Value = "  +0.19%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'5.525"
$ws.Range("E13").Value = "  +0.15%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "'6.744"
$ws.Range("E14").Value = "  +0.69%  "

# Row 15 - TRON
$ws.Range("D15").Value = "'0.07208"
$ws.Range("E15").Value = "  +0.05%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'86.06"
$ws.Range("E16").Value = "  +4.09%  "

# Row 17 - BinanceUSD
$ws.Range("D17").Value = "'1.038"
$ws.Range("E17").Value = "  +0.34%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.000009102"
$ws.Range("E18").Value = "  +0.68%  "

# Row 19 - Dai
$ws.Range("D19").Value = "'1.032"
$ws.Range("E19").Value = "  +0.26%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "'15.56"
$ws.Range("E20").Value = "  +0.67%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "27.771.38"
$ws.Range("E21").Value = "  +0.69%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.300"

# Row 23 - Cosmos
$ws.Range("D23").Value = "'11.27"
$ws.Range("E23").Value = "  +0.29%  "

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "2.090.13"
$ws.Range("E24").Value = "  +0.62%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'2.073"
$ws.Range("E25").Value = "  +6.45%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'159.08"
$ws.Range("E26").Value = "  +0.83%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'18.76"
$ws.Range("E27").Value = "  +0.05%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "'5.359"
$ws.Range("E28").Value = "  +1.02%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'1.989"
$ws.Range("E29").Value = "  +2.53%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "'118.68"
$ws.Range("E30").Value = "  +1.90%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.09117"
$ws.Range("E31").Value = "  +0.64%  "

# Row 32 - was ARBITRUM, now ImmutableX
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.7732"
$ws.Range("E32").Value = "  +0.61%  "

# Row 33 - was ImmutableX, now ARBITRUM
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.212"
$ws.Range("E33").Value = "  +0.16%  "

# Row 34 - HuobiToken
$ws.Range("D34").Value = "'3.044"
$ws.Range("E34").Value = "  +5.25%  "

# Row 35 - Filecoin
$ws.Range("D35").Value = "'4.599"
$ws.Range("E35").Value = "  +1.47%  "

# Row 36 - Frax
$ws.Range("E36").Value = "  +0.15%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "'1.154"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.01986"
$ws.Range("E38").Value = "  +0.37%  "

# Row 39 - Hedera
$ws.Range("D39").Value = "'0.05317"
$ws.Range("E39").Value = "  +0.52%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "'2.862"
$ws.Range("E40").Value = "  +1.34%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "'0.5216"
$ws.Range("E41").Value = "  +0.77%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "'6.957"
$ws.Range("E42").Value = "  +3.24%  "

# Row 43 - Algorand
$ws.Range("D43").Value = "'0.1678"
$ws.Range("E43").Value = "  +0.24%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "'8.783"
$ws.Range("E44").Value = "  +2.16%  "

# Row 45 - was Quant, now EnergySwap
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'10.83"
$ws.Range("E45").Value = "  +1.90%  "

# Row 46 - was EnergySwap, now Quant
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'110.29"
$ws.Range("E46").Value = "  +1.22%  "

# Row 47 - PaxDollar
$ws.Range("E47").Value = "  +0.12%  "

# Row 48 - Cronos
$ws.Range("D48").Value = "'0.06559"
$ws.Range("E48").Value = "  +2.42%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "'1.714"
$ws.Range("E49").Value = "  +0.02%  "

# Row 50 - Decentraland
$ws.Range("D50").Value = "'0.4729"
$ws.Range("E50").Value = "  +1.41%  "

# Row 51 - RenderToken
$ws.Range("D51").Value = "'1.887"
$ws.Range("E51").Value = "  -0.37%  "
